# Circle Language Spec Plan: Set font to Calibri for non-heading text.
#
# 1. Change the "Normal" style's font to Calibri, 11pt (sz=22 half-points),
#    which cascades to every paragraph/heading that doesn't override its
#    own rFonts/sz.
# 2. Move the "_GoBack" bookmark (Word's "last edit position" marker) from
#    wherever it was left in the old revision to the point in the text
#    where the author's cursor ended up: right after "...looking up",
#    before " how exactly something was done." This naturally splits the
#    run in two, matching the target markup.

$d = $word.ActiveDocument

# --- 1. Font change: Normal style -> Calibri 11pt ---------------------
$normal = $d.Styles("Normal")
$normal.Font.Name = "Calibri"
$normal.Font.Size = 11

# --- 2. Reposition the _GoBack bookmark --------------------------------
$find = $d.Content
$found = $find.Find.Execute("for looking up", $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0)
if ($found) {
    $splitPoint = $find.End
    $goBack = $d.Range($splitPoint, $splitPoint)
    # Adding a bookmark with a name that already exists moves it (Word
    # only allows one bookmark per name), so this both removes the old
    # "_GoBack" (which used to wrap almost the whole rest of the document)
    # and creates the new, empty one at the cursor's last position.
    $d.Bookmarks.Add("_GoBack", $goBack)
}
